$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Nombre" column (A) for existing rows 10-12, which was left blank ---
$ws.Range("A10").Value = "Federico Speroni"
$ws.Range("A11").Value = "Federico Speroni"
$ws.Range("A12").Value = "Federico Speroni"

# --- New row 13 ---
$ws.Range("A13").Value = "Federico Speroni"
$ws.Range("B13").Value = 42837
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = "Creación de Anteproyecto"
$ws.Range("E13").Value = "Plan de SQA, estandares definidos y convenciones"

# --- New row 14 ---
$ws.Range("B14").Value = 42838
$ws.Range("B12").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = "Creación de Anteproyecto"
$ws.Range("E14").Value = "Plan de testing, Pland de SCM, Plan de capacitación. Documentación de EncuestaUsuarios1"

# --- Update the visible selection/scroll position like the author left it ---
$ws.Range("E16").Select()
